$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("F7").Value = "color = ""red""`nside_1 = 12`nside_2 = 12`nside_3 = 12"
$ws.Range("G7").Value = "Object created with expected attribute values"

# Row 8
$ws.Range("F8").Value = "color = """"`nside_1 = 12`nside_2 = 12`nside_3 = 12"
$ws.Range("G8").Value = "ValueError"

# Row 9
$ws.Range("F9").Value = "color = ""red""`nside_1 = ""hi""`nside_2 = 12`nside_3 = 12"
$ws.Range("G9").Value = "ValueError"

# Row 10
$ws.Range("F10").Value = "color = ""red""`nside_1 = 12`nside_2 = ""hi""`nside_3 = 12"
$ws.Range("G10").Value = "ValueError"

# Row 11
$ws.Range("F11").Value = "color = ""red""`nside_1 = 12`nside_2 = 12`nside_3 = ""hi"""
$ws.Range("G11").Value = "ValueError"

# Row 12
$ws.Range("E12").Value = "triangle = Triangle(""Red"", 12, 12, 12)"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "String returned in correct format"

# Row 13
$ws.Range("E13").Value = "triangle = Triangle(""Red"", 12, 12, 12)"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "triangle.area = 62.35"

# Row 14
$ws.Range("E14").Value = "triangle = Triangle(""Red"", 12, 12, 12)"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "triangle.area = 36"

# Sheet view: scroll position (best effort) + active selection
try { $excel.ActiveWindow.ScrollRow = 7 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 4 } catch {}
$ws.Range("G14").Select()
